$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the roster with the new names / emails (existing rows are overwritten
# in top-to-bottom order so shared-string indices line up the way Excel would
# naturally allocate them).
$ws.Range("B2").Value = "Darshan"
$ws.Range("C2").Value = "pasne.d@husky.neu.edu"

$ws.Range("B3").Value = "Saman"
$ws.Range("C3").Value = "sood.s@husky.neu.edu"

$ws.Range("B4").Value = "Shail"
$ws.Range("C4").Value = "shail@ccs.neu.edu"

$ws.Range("B5").Value = "Vaibhav"
$ws.Range("C5").Value = "dave.v@husky.neu.edu"

$ws.Range("B6").Value = "John"
$ws.Range("C6").Value = "snow.j@husky.neu.edu"

$ws.Range("B7").Value = "Danny"
$ws.Range("C7").Value = "danny.d@husky.neu.edu"

$ws.Range("B8").Value = "Erica"
$ws.Range("C8").Value = "sniper.e@husky.neu.edu"

$ws.Range("B9").Value = "Flurry"
$ws.Range("C9").Value = "majin.f@husky.neu.edu"

$ws.Range("B10").Value = "Gara"
$ws.Range("C10").Value = "hawking.g@husky.neu.edu"

$ws.Range("B11").Value = "Max"
$ws.Range("C11").Value = "max@x.com"

$ws.Range("B12").Value = "Kat"
$ws.Range("C12").Value = "kat@x.com"

# The old mailto: hyperlinks no longer apply to the refreshed e-mail addresses.
$ws.Hyperlinks.Delete()

# Leave the selection where the author ended up after entering the new data.
[void]$ws.Range("C19").Select()
